$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

for ($r = 2; $r -le 62; $r++) {
    $cell = $ws.Cells.Item($r, 4)  # Column D
    $val = $cell.Value2
    if ($val -eq "MOLLY MCNINCH") {
        $cell.Value2 = "T"
    } elseif ($val -eq "STUDENT") {
        $cell.Value2 = "S"
    }
}
